$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "userid" header in F1, copying the header style (bold, border, centered)
# from the neighboring header cell E1.
$ws.Range("F1").Value = "userid"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update lastlogin timestamps for existing rows
$ws.Range("E2").Value = "2020-03-02 19:57:05.57S"
$ws.Range("E3").Value = "2020-03-02 17:53:19.53S"

# Add new userid numeric values for the two rows
$ws.Range("F2").Value = 100
$ws.Range("F3").Value = 101
